# "files with surveyor names"
# Insert a new "surveyor" column before the existing "notes" column (D),
# pushing "notes" to column E, then fill in the surveyor name for every
# data row (Grace for the first and third block of rows, Erik for the
# second and fourth block).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shift the existing "notes" column (D) right to make room for the new
# "surveyor" column.
$ws.Columns("D").Insert()

# New header for column D.
$ws.Range("D1").Value2 = "surveyor"

# Surveyor name per data row (2-23).
$surveyors = @{
    2  = "Grace"
    3  = "Grace"
    4  = "Grace"
    5  = "Grace"
    6  = "Grace"
    7  = "Grace"
    8  = "Erik"
    9  = "Erik"
    10 = "Erik"
    11 = "Erik"
    12 = "Erik"
    13 = "Grace"
    14 = "Grace"
    15 = "Grace"
    16 = "Grace"
    17 = "Grace"
    18 = "Grace"
    19 = "Erik"
    20 = "Erik"
    21 = "Erik"
    22 = "Erik"
    23 = "Erik"
}

foreach ($row in $surveyors.Keys) {
    $ws.Cells.Item($row, 4).Value2 = $surveyors[$row]
}

# Match the author's recorded selection after the edit.
$ws.Range("F24").Select()
